$d = $word.ActiveDocument

# --- Fix 1: footnote w:id="29" (the ninth footnote) currently ends with a
# stray "aa" that shouldn't be there ("སཱཙྪ། སྣར་ཐང་། པེ་ཅིན།aa" ->
# "སཱཙྪ། སྣར་ཐང་། པེ་ཅིན།"). Rewrite the footnote's text without the "aa".
$fnAA = $d.Footnotes.Item(9)
$fnAA.Range.Text = "སཱཙྪ། སྣར་ཐང་། པེ་ཅིན།"

# --- Fix 2: footnote w:id="45" (the last footnote, reference at the very
# end of the body paragraph) is an empty/blank note (just "།") that was
# added by mistake. Remove the footnote entirely -- this deletes both the
# footnoteReference run in the body and the footnote definition itself.
$fnEmpty = $d.Footnotes.Item(25)
$fnEmpty.Delete()
